# Update generated_model_inputs.xlsx:
#  1. Refresh the run timestamps logged on the "Log" sheet (column A, rows 2-33)
#     to reflect the latest run of the input-generation process.
#  2. Rename the "FA_STORE" references that should point at the Gladstone
#     export silo to "FA_EXPSILO_STORE" on the "SHIP_ROUTES" sheet
#     (Product 1 Store for routes 6/7/8, and Product 2 Store for route 4)
#     so Fly Ash shipments resolve to the correct storage location.

$wb = $excel.ActiveWorkbook

# --- 1. Log sheet timestamps -------------------------------------------------
$logSheet = $wb.Worksheets.Item("Log")

$timestamps = @{
    2  = 46008.50060285225
    3  = 46008.50060357437
    4  = 46008.5006038809
    5  = 46008.50060424916
    6  = 46008.50060460019
    7  = 46008.50060492102
    8  = 46008.50060518074
    9  = 46008.50060554549
    10 = 46008.5006058962
    11 = 46008.50060618638
    12 = 46008.50060648269
    13 = 46008.50060678002
    14 = 46008.50060706501
    15 = 46008.50060739729
    16 = 46008.50060741292
    17 = 46008.50060744611
    18 = 46008.50060751093
    19 = 46008.50060753969
    20 = 46008.50060760233
    21 = 46008.50060766236
    22 = 46008.50060770309
    23 = 46008.50060781257
    24 = 46008.50060784438
    25 = 46008.50061345316
    26 = 46008.50061783077
    27 = 46008.50062267205
    28 = 46008.50062741721
    29 = 46008.50063151732
    30 = 46008.50063740841
    31 = 46008.50064160447
    32 = 46008.50064491263
    33 = 46008.5006483852
}

foreach ($row in $timestamps.Keys) {
    $logSheet.Cells.Item($row, 1).Value = $timestamps[$row]
}

# --- 2. SHIP_ROUTES FA_STORE -> FA_EXPSILO_STORE -----------------------------
$routes = $wb.Worksheets.Item("SHIP_ROUTES")

$routes.Cells.Item(6, 7).Value = "FA_EXPSILO_STORE"   # G6 - Route 6, Product 1 Store
$routes.Cells.Item(6, 8).Value = "FA_EXPSILO_STORE"   # H6 - Route 7, Product 1 Store
$routes.Cells.Item(6, 9).Value = "FA_EXPSILO_STORE"   # I6 - Route 8, Product 1 Store
$routes.Cells.Item(8, 5).Value = "FA_EXPSILO_STORE"   # E8 - Route 4, Product 2 Store
